$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "ingredients"

# Update the Photo URL for the Hot Dog row (E23)
$ws.Range("E23").Value = "http://farm9.staticflickr.com/8511/8598717582_2d4ec7e7e4_z.jpg"

# Scroll the view so row 5 is at the top, then update the selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E26").Select()
